$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.351.95'
$ws.Range("E2").Value = '  -0.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.688.88'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '678.51'
$ws.Range("E5").Value = '  -1.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.15'
$ws.Range("E6").Value = '  -1.69%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -0.36%  '

$ws.Range("E9").Value = '  -0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.17'
$ws.Range("E10").Value = '  -2.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.442'
$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000232'
$ws.Range("E12").Value = '  -1.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.308.75'
$ws.Range("E13").Value = '  -0.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.35'
$ws.Range("E14").Value = '  -1.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.680.73'
$ws.Range("E15").Value = '  -0.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.361.26'
$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("E17").Value = '  +2.87%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.01'
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.47'
$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.37'
$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("E21").Value = '  -0.71%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.652'
$ws.Range("E22").Value = '  -0.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.00'
$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.834.35'
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  -4.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.90'
$ws.Range("E27").Value = '  -2.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.15'
$ws.Range("E28").Value = '  -1.09%  '

$ws.Range("E29").Value = '  -0.44%  '

$ws.Range("E30").Value = '  -3.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.56'
$ws.Range("E31").Value = '  -3.03%  '

$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("E33").Value = '  -2.63%  '

$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.678.20'
$ws.Range("E35").Value = '  +0.64%  '

$ws.Range("E36").Value = '  -2.84%  '

$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.29'
$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.25'
$ws.Range("E40").Value = '  -2.67%  '

$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '169.88'
$ws.Range("E43").Value = '  +3.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.942'
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.11'
$ws.Range("E45").Value = '  -2.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.13'
$ws.Range("E46").Value = '  -6.86%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.70'
$ws.Range("E47").Value = '  -1.63%  '

$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000278'
$ws.Range("E48").Value = '  -0.78%  '

$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.11'
$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("E50").Value = '  -2.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.83'
$ws.Range("E51").Value = '  -2.08%  '
